$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values are written as text (not auto-converted to numbers)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.904.41'
$ws.Range("E2").Value = '  -5.57%  '
$ws.Range("D3").Value = '1.821.96'
$ws.Range("E3").Value = '  -4.29%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.60%  '
$ws.Range("D5").Value = '327.60'
$ws.Range("E5").Value = '  -3.02%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.42%  '
$ws.Range("D7").Value = '0.4629'
$ws.Range("E7").Value = '  -2.92%  '
$ws.Range("D8").Value = '0.3843'
$ws.Range("E8").Value = '  -3.88%  '
$ws.Range("D9").Value = '45.76'
$ws.Range("E9").Value = '  -3.18%  '
$ws.Range("D10").Value = '0.07843'
$ws.Range("E10").Value = '  -2.60%  '
$ws.Range("D11").Value = '0.9589'
$ws.Range("E11").Value = '  -3.35%  '
$ws.Range("D12").Value = '21.84'
$ws.Range("E12").Value = '  -5.85%  '
$ws.Range("D13").Value = '1.832.81'
$ws.Range("E13").Value = '  -2.59%  '
$ws.Range("D14").Value = '5.643'
$ws.Range("E14").Value = '  -4.81%  '
$ws.Range("D15").Value = '6.848'
$ws.Range("E15").Value = '  -3.79%  '
$ws.Range("D16").Value = '0.06870'
$ws.Range("E16").Value = '  +0.82%  '
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  -0.58%  '
$ws.Range("D18").Value = '86.46'
$ws.Range("E18").Value = '  -3.07%  '
$ws.Range("D19").Value = '0.000009928'
$ws.Range("E19").Value = '  -2.91%  '
$ws.Range("D20").Value = '16.64'
$ws.Range("E20").Value = '  -4.07%  '
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("D22").Value = '27.911.90'
$ws.Range("E22").Value = '  -5.51%  '
$ws.Range("E23").Value = '  -3.79%  '
$ws.Range("D24").Value = '10.95'
$ws.Range("E24").Value = '  -6.26%  '
$ws.Range("D25").Value = '2.099'
$ws.Range("E25").Value = '  -2.73%  '
$ws.Range("D26").Value = '2.028.07'
$ws.Range("E26").Value = '  -3.54%  '
$ws.Range("D27").Value = '152.00'
$ws.Range("E27").Value = '  -3.07%  '
$ws.Range("D28").Value = '19.13'
$ws.Range("E28").Value = '  -2.43%  '
$ws.Range("D29").Value = '5.767'
$ws.Range("E29").Value = '  -11.87%  '
$ws.Range("D30").Value = '1.967'
$ws.Range("E30").Value = '  -4.35%  '
$ws.Range("D31").Value = '116.48'
$ws.Range("E31").Value = '  -2.30%  '
$ws.Range("D32").Value = '0.9334'
$ws.Range("E32").Value = '  -6.37%  '
$ws.Range("D33").Value = '0.09224'
$ws.Range("E33").Value = '  -3.33%  '
$ws.Range("D34").Value = '5.272'
$ws.Range("E34").Value = '  -3.59%  '
$ws.Range("D35").Value = '1.314'
$ws.Range("E35").Value = '  -5.38%  '
$ws.Range("D36").Value = '3.343'
$ws.Range("E36").Value = '  -5.46%  '
$ws.Range("D37").Value = '0.05927'
$ws.Range("E37").Value = '  -8.50%  '
$ws.Range("D38").Value = '0.02142'
$ws.Range("E38").Value = '  -4.44%  '
$ws.Range("D39").Value = '1.142'
$ws.Range("E39").Value = '  -4.30%  '
$ws.Range("E40").Value = '  -0.35%  '
$ws.Range("D41").Value = '7.532'
$ws.Range("E41").Value = '  -2.45%  '
$ws.Range("D42").Value = '0.5565'
$ws.Range("E42").Value = '  -4.41%  '
$ws.Range("D43").Value = '9.921'
$ws.Range("E43").Value = '  -5.81%  '
$ws.Range("E44").Value = '  -3.27%  '
$ws.Range("D45").Value = '1.232'
$ws.Range("E45").Value = '  -2.77%  '
$ws.Range("D46").Value = '2.209'
$ws.Range("E46").Value = '  -10.14%  '
$ws.Range("D47").Value = '11.53'
$ws.Range("E47").Value = '  -5.32%  '
$ws.Range("D48").Value = '0.5234'
$ws.Range("E48").Value = '  -4.52%  '
$ws.Range("D49").Value = '0.06985'
$ws.Range("E49").Value = '  -5.70%  '
$ws.Range("D50").Value = '1.817'
$ws.Range("E50").Value = '  -7.14%  '
$ws.Range("D51").Value = '111.91'
$ws.Range("E51").Value = '  -3.47%  '

# Restore default style (remove explicit text number format) to match original formatting
$ws.Range("D2:D51").Style = "Normal"
